$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column AC: date header "2020-04-02" -----------------------------
# Copy the format (style, border, bold, centering) from the previous
# header cell AB1 so AC1 matches the rest of the header row.
$ws.Range("AB1").Copy()
$ws.Range("AC1").PasteSpecial(-4122)
# Build the date label from existing text (AB1 = "2020-04-01") so the
# engine treats it as a plain string instead of auto-converting a typed
# "2020-04-02" literal into a date serial number.
$ws.Range("AC1").Formula = "=LEFT(AB1,9)&""2"""
$ws.Range("AC1").Copy()
$ws.Range("AC1").PasteSpecial(-4163)

# --- New column AC: cumulative counts for rows 2-22 -----------------------
$acValues = @(2, 13, 14, 10, 1, 15, 4, 0, 20, 29, 247, 46, 21, 4, 8, 6, 13, 16, 35, 79, 584)
for ($i = 0; $i -lt $acValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 29).Value = $acValues[$i]
}

# --- New row 23: "Region Gotland" ------------------------------------------
# Materialize the blank cells B23:AB23 (present but empty, like B2/C2 in
# row 2) by pasting formats from the corresponding cells of row 2.
$ws.Range("B2:AB2").Copy()
$ws.Range("B23:AB23").PasteSpecial(-4122)

# Region label in column A, styled like the other region names.
$ws.Range("A2").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Region Gotland"

# First reported value for Region Gotland, in the new AC column.
$ws.Cells.Item(23, 29).Value = 1
